$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# The "Language of Interpreter: ______." run and the following
# "{% endif %}" run get merged into a single run/text node:
#   ": __________________________."  +  "{% endif %}"
#   -> ": __________________________.{% endif %}"
$found1 = $d.Content.Find.Execute(
    ": __________________________.{% endif %}", $true, $false, $false,
    $false, $false, $true, 1, $false,
    ": __________________________.{% endif %}", 2)
Write-Output "change1 found=$found1"

# --- Change 2 -------------------------------------------------------------
# Append a new sentence about the Victim's Attorney right after the
# existing "...County Jail: PS   EM;" text in the COS paragraph.
$rng = $d.Content
$found2 = $rng.Find.Execute(
    "County Jail: PS   EM;", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
Write-Output "change2 found=$found2"

$rng.Collapse(0)
$rng.InsertAfter(" ")
$rng.Collapse(0)
$rng.InsertAfter("Victim’s Attorney (if applicable): PS   OS   EM")

Write-Output "done"
